$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 20: it used to hold the "Big Cities have 4 tanks..." note (now removed
# from the TODO list). It now takes over the content that row 21 used to
# hold ("Apparently "Anthrakia" is registered as a big city"), and its
# formatting is promoted from the alternating "green" style to the "plain"
# style that row 21 used to have - copy that formatting from row 21 before
# row 21 gets overwritten with new content below.
# ---------------------------------------------------------------------------
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "Apparently ""Anthrakia"" is registered as a big city"

# ---------------------------------------------------------------------------
# Row 21: becomes a brand-new entry (date/time/description), keeping the
# "plain" style it already has.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = 41956
$ws.Range("B21").Value = "18.40"
$ws.Range("C21").Value = "Recheck the R3F-Config-Lists in plugins\R3F_LOG\addons_config\A3_vanilla_1.22"

# ---------------------------------------------------------------------------
# Row 22 (new): same date/time as row 21 ("done" item), rendered with a new
# strikethrough dark-green font to mark it as resolved.
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = 41956
$ws.Range("B22").Value = "18.40"
$ws.Range("C22").Value = "Recode the number of patrols, for each target"

$ws.Range("A22").Font.Color = 5287936
$ws.Range("A22").Font.Strikethrough = $true
$ws.Range("A22").NumberFormat = "mm-dd-yy"

$ws.Range("B22").Font.Color = 5287936
$ws.Range("B22").Font.Strikethrough = $true
$ws.Range("B22").NumberFormat = "@"

$ws.Range("C22").Font.Color = 5287936
$ws.Range("C22").Font.Strikethrough = $true

# ---------------------------------------------------------------------------
# Rows 23-24 (new): more items for the same date/time, normal "plain" style
# (same formatting as row 21 - copy it across).
# ---------------------------------------------------------------------------
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A23:C23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = 41956
$ws.Range("B23").Value = "18.40"
$ws.Range("C23").Value = "Some ammoboxes exploded, they seem to take damage, when manually setting allowdamage to false (by console) the don't take damage anymore"

$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A24:C24").PasteSpecial(-4122) | Out-Null
$ws.Range("A24").Value = 41956
$ws.Range("B24").Value = "18.40"
$ws.Range("C24").Value = "Communication options disappear when a save mission is loaded"

# ---------------------------------------------------------------------------
# Rows 25-26 (new): a new date/time group, using the "first row of a group"
# style (same formatting as row 15 - copy it across).
# ---------------------------------------------------------------------------
$ws.Range("A15:C15").Copy() | Out-Null
$ws.Range("A25:C25").PasteSpecial(-4122) | Out-Null
$ws.Range("A25").Value = 41957
$ws.Range("B25").Value = "21.00"
$ws.Range("C25").Value = "Logical mistake when determining the amount of enemies per target, the whole determination has been changed now"

$ws.Range("A15:C15").Copy() | Out-Null
$ws.Range("A26:C26").PasteSpecial(-4122) | Out-Null
$ws.Range("A26").Value = 41957
$ws.Range("B26").Value = "21.00"
$ws.Range("C26").Value = "Recoded the amount of loot, god tier weapons were almost impossible to get"

# ---------------------------------------------------------------------------
# Move the selection like the recorded session did (last cell now C27, one
# row below the new last data row).
# ---------------------------------------------------------------------------
$ws.Range("C27").Select() | Out-Null

Write-Output "done"
